$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stale "_GoBack" bookmark that currently sits on the empty
#    first paragraph of the document (left over from a previous edit).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. The title-page date run currently reads "XX 2019" and is wrapped by the
#    "_Toc322009854" / "_Toc322027047" bookmarks (spanning the whole "XX 2019"
#    text). The edit bumps the year to "2020"; the last two keystrokes of
#    that edit happened after "XX 20", so the TOC bookmarks now close right
#    after "XX 20" instead of at the end of the run, and the new text ("20")
#    that follows starts outside of them.
# ---------------------------------------------------------------------------
$tocBm1 = $d.Bookmarks("_Toc322009854")
$tocBm2 = $d.Bookmarks("_Toc322027047")
$tocStart1 = $tocBm1.Start
$tocStart2 = $tocBm2.Start
$tocEnd = $tocBm1.End

# The last two characters of the run ("19") are where the new "20" goes;
# the TOC bookmarks will close right before those two characters.
$splitPoint = $tocEnd - 2

$tocBm1.Delete()
$tocBm2.Delete()
$d.Bookmarks.Add("_Toc322009854", $d.Range($tocStart1, $splitPoint))
$d.Bookmarks.Add("_Toc322027047", $d.Range($tocStart2, $splitPoint))

# Replace the trailing "19" with "20", turning "XX 2019" into "XX 2020".
$tailRange = $d.Range($splitPoint, $tocEnd)
$tailRange.Text = "20"

# ---------------------------------------------------------------------------
# 3. Word drops a fresh "_GoBack" bookmark at the spot where the edit
#    finished, i.e. right after the newly typed "20" and before the
#    following page break.
# ---------------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $d.Range($tocEnd, $tocEnd))
